$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns EB..EO (before deletion) hold the comments that must shift one
# column to the left (into EA..EN) once the duplicate EA column is removed.
$srcCols = @("EB","EC","ED","EE","EF","EG","EH","EI","EJ","EK","EL","EM","EN","EO")
$dstCols = @("EA","EB","EC","ED","EE","EF","EG","EH","EI","EJ","EK","EL","EM","EN")

# Capture the comment text before touching the sheet.
$texts = @{}
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $ref = $srcCols[$i] + "4"
    $cmt = $ws.Range($ref).Comment
    if ($cmt -ne $null) {
        $texts[$dstCols[$i]] = $cmt.Text()
    }
}

# Remove the duplicate "NOMBRE_ARTICULO_2" column (EA); this shifts every
# subsequent column (data, shared-string refs, dimension) one slot left.
[void]$ws.Range("EA1:EA4").EntireColumn.Delete()

# The comments collection is NOT re-anchored by the column delete, so move
# each caption manually into its new (post-delete) cell.
for ($i = 0; $i -lt $dstCols.Length; $i++) {
    $dstRef = $dstCols[$i] + "4"
    $newText = $texts[$dstCols[$i]]
    $cmt = $ws.Range($dstRef).Comment
    if ($cmt -ne $null -and $newText -ne $null) {
        [void]$cmt.Text($newText)
    }
}

# EO4's comment is now orphaned (its text already lives on EN4); drop it.
$tailComment = $ws.Range("EO4").Comment
if ($tailComment -ne $null) {
    [void]$tailComment.Delete()
}
